$wb = $excel.ActiveWorkbook

# --- NK sheet: update the percentage factors that drive the JRP criteria table ---
$nk = $wb.Worksheets.Item("NK")
$nk.Range("F14:F19").Value = 1.19
$nk.Range("U14:U19").Value = 1.14
$nk.Range("AA14:AA19").Value = 1.21

# --- JRP sheet: update the criteria row labels (percentages changed) ---
$jrp = $wb.Worksheets.Item("JRP")
$jrp.Range("A4:A6").Value = "NK2 +19% Mäner"
$jrp.Range("A7:A9").Value = "NK2 +18% Frauen"
$jrp.Range("A10:A12").Value = "NK1 +14% Mäner"
$jrp.Range("A13:A15").Value = "NK1 +21% Frauen"

# --- Restore / set per-sheet selection state ---
$nk.Activate()
$nk.Range("AA14").Select()

# --- Make JRP the active sheet/tab, with its previous selection ---
$jrp.Activate()
$jrp.Range("A15").Select()
